# Add a new worksheet "PQ CAMION" at the end of the workbook, with the
# same look & feel (styles) as the neighbouring "PQ AUSEN" sheet, and
# populate it with the technician/truck-scan data.

$wb = $excel.ActiveWorkbook

$styleSource = $wb.Worksheets.Item("PQ AUSEN")
$lastSheet   = $wb.Worksheets.Item($wb.Worksheets.Count)

$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "PQ CAMION"

# ---- Header row -----------------------------------------------------
$ws.Range("A1").Value = "FECHA SCAN PET"
$ws.Range("B1").Value = "CEDULA TECNICO"
$ws.Range("C1").Value = "NOMBRE TECNICO "
$ws.Range("D1").Value = "Centro de costos"
$ws.Range("E1").Value = "OBSERVACIÓN"
$ws.Range("F1").Value = "Cargo colaborador"

# Reuse the bold/bordered header style that is already used on "PQ AUSEN".
$styleSource.Range("A1:E1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)

# ---- Data rows --------------------------------------------------------
$centro = "1197001 - Aprovisionamiento FTTH C-CW5903-23 BOG"
$cargo  = "TECNICO II / INSTALADOR FTTH"

$fechas   = 45414, 45414, 45415, 45416, 45416, 45418, 45419
$cedulas  = 1022355666, 1006886735, 1022355666, 1006886735, 1022355666, 1022355666, 1022355666
$nombres  = "JUAN DAVID BECERRA PEÑA", "WILBERTO DARIO PEREZ TILLER", "JUAN DAVID BECERRA PEÑA", "WILBERTO DARIO PEREZ TILLER", "JUAN DAVID BECERRA PEÑA", "JUAN DAVID BECERRA PEÑA", "JUAN DAVID BECERRA PEÑA"

for ($i = 0; $i -lt $fechas.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $fechas[$i]
    $ws.Cells.Item($row, 2).Value = $cedulas[$i]
    $ws.Cells.Item($row, 3).Value = $nombres[$i]
    $ws.Cells.Item($row, 4).Value = $centro
    $ws.Cells.Item($row, 6).Value = $cargo
}

# Reuse the date-formatted style from "PQ AUSEN" column A for the date column.
$styleSource.Range("A2:A8").Copy()
$ws.Range("A2:A8").PasteSpecial(-4122)

# Re-apply the values after the paste-special (paste-special(formats) does
# not touch values, but make sure dates are numbers, not re-evaluated).
for ($i = 0; $i -lt $fechas.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $fechas[$i]
}

$ws.Range("A1").Select()
